# Update last_edited_time values and discount-rate totals for the "Thang 7" rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D3, D4, D5, D7, D13 -> new last_edited_time value (stored as plain text, same as source)
$ws.Range("D3").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D4").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D5").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D7").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D13").Value = "2024-07-25T15:02:00.000Z"

# Update numeric totals on row 5 (Chi tieu / Luy ke formula numbers)
$ws.Range("W5").Value = 30374000
$ws.Range("AA5").Value = 15576000
